$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates
# Leading apostrophe forces Excel to keep these as plain text (matching the
# original inlineStr/text cell type) instead of auto-converting
# number-looking values (e.g. "253.03") into numeric values.
$ws.Range("D2").Value = "'35.129.96"
$ws.Range("D3").Value = "'1.901.82"
$ws.Range("D5").Value = "'253.03"
$ws.Range("D8").Value = "'41.42"
$ws.Range("D10").Value = "'52.69"
$ws.Range("D12").Value = "'0.0982"
$ws.Range("D13").Value = "'2.176.73"
$ws.Range("D14").Value = "'13.00"
$ws.Range("D17").Value = "'1.881.27"
$ws.Range("D18").Value = "'35.139.40"
$ws.Range("D19").Value = "'73.44"
$ws.Range("D20").Value = "'0.0₃0831"
$ws.Range("D21").Value = "'243.42"
$ws.Range("D22").Value = "'12.92"
$ws.Range("D26").Value = "'2.28"
$ws.Range("D27").Value = "'166.56"
$ws.Range("D29").Value = "'18.46"
$ws.Range("D31").Value = "'4.128.83"
$ws.Range("D35").Value = "'1.57"
$ws.Range("D38").Value = "'0.850"
$ws.Range("D39").Value = "'2.00"
$ws.Range("D40").Value = "'102.40"
$ws.Range("D41").Value = "'17.24"
$ws.Range("D44").Value = "'0.0649"
$ws.Range("D45").Value = "'1.318.95"
$ws.Range("D48").Value = "'2.74"
$ws.Range("D49").Value = "'12.27"
$ws.Range("D50").Value = "'6.57"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  +4.61%  "
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  +5.00%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +12.89%  "
$ws.Range("E33").Value = "  +6.84%  "
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  -7.60%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  +14.46%  "
$ws.Range("E41").Value = "  +7.86%  "
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("E51").Value = "  +5.59%  "
